$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39.78, 0.16, 23.74857687950134, 16.28125),
    @(39.78, 0.16, 38.82865643501282, 18.640625),
    @(39.78, 0.16, 91.12993431091309, 22.765625),
    @(39.06, 0.04000000000000001, 10.32614302635193, 5.03125)
)

$startRow = 220
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
